$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook stores every data value (including numeric-looking prices)
# as TEXT (inline/shared strings), matching the crypto-tracker source data
# (trailing zeros like "0.530" and "12.20" must be preserved literally).
# For column D (Price) we therefore force Text format before assigning the
# value so Excel does not silently coerce it to a Number (which would drop
# trailing zeros / renormalize the text); the style is then reset back to
# "Normal" so no stray NumberFormat is left on the cell.

# Row 2 (Bitcoin)
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "43.454.64"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = "  +2.59%  "

# Row 3 (Ethereum)
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.313.94"
$ws.Cells.Item(3, 4).Style = "Normal"

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 (BNB)
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "311.18"
$ws.Cells.Item(5, 4).Style = "Normal"

# Row 6 (Solana)
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "104.08"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  +6.46%  "

# Row 7 (XRP)
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.534"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  +1.12%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 (Cardano)
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.530"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  +8.25%  "

# Row 10 (Avalanche)
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "36.72"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  +4.01%  "

# Row 11 (OKB)
$ws.Range("E11").Value = "  +1.29%  "

# Row 12 (Dogecoin)
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0814"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "  +0.42%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  -1.26%  "

# Row 14 (Polkadot)
$ws.Range("E14").Value = "  +2.52%  "

# Row 15 (WrappedliquidstakedEther2.0)
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "2.670.93"
$ws.Cells.Item(15, 4).Style = "Normal"

# Row 16 (Chainlink)
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "15.13"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  +3.28%  "

# Row 17 (WrappedEther)
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.310.22"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  +2.09%  "

# Row 18 (Polygon)
$ws.Range("E18").Value = "  +2.42%  "

# Row 19 (WrappedBTC)
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "43.362.31"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +2.71%  "

# Row 20 (InternetComputer(DFINITY))
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "12.18"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "

# Row 21 (ShibaInu)
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.0₃0927"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  +2.26%  "

# Row 22 (Uniswap)
$ws.Range("E22").Value = "  +3.36%  "

# Row 23 (Litecoin)
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "68.17"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = "  +0.77%  "

# Row 24 (BitcoinCash)
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "242.82"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = "  +2.56%  "

# Row 25 (ImmutableX)
$ws.Range("E25").Value = "  +2.57%  "

# Row 26 (PancakeSwap)
$ws.Range("E26").Value = "  +0.95%  "

# Row 27 (Dai)
$ws.Range("E27").Value = "  +0.17%  "

# Row 28 (EthereumClassic)
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "24.92"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  +5.56%  "

# Row 29 (Toncoin)
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.35"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  +10.34%  "

# Row 30 (InjectiveProtocol)
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "37.14"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  -1.10%  "

# Row 31 (Cosmos)
$ws.Range("E31").Value = "  +0.54%  "

# Row 32 (Monero)
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "167.87"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  +2.50%  "

# Row 34 (FirstDigitalUSD)
$ws.Range("E34").Value = "  +0.00%  "

# Row 35 (Celestia)
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "18.41"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +3.67%  "

# Row 36 (WEMIXToken)
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.53"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  +6.71%  "

# Row 37 (Hedera)
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0745"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "  +1.17%  "

# Row 38 (LidoDAOToken)
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.07"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "

# Row 39 (ARBITRUM)
$ws.Range("E39").Value = "  +3.19%  "

# Row 40 (Kaspa -> RenderToken)
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "4.47"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  +7.05%  "

# Row 41 (RenderToken -> Kaspa)
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.106"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  +1.92%  "

# Row 42 (Stellar)
$ws.Range("E42").Value = "  +0.69%  "

# Row 43 (ApeXProtocol)
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.75"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  +21.31%  "

# Row 44 (VeChain)
$ws.Range("E44").Value = "  +3.62%  "

# Row 45 (Maker)
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.992.33"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  +2.23%  "

# Row 46 (EnergySwap)
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "19.09"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  +1.63%  "

# Row 47 (NEARProtocol)
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.07"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  +3.50%  "

# Row 48 (FraxShare)
$ws.Range("E48").Value = "  +2.07%  "

# Row 49 (MultiversX)
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "56.22"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  +3.83%  "

# Row 50 (HuobiToken)
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.96"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "  +2.69%  "

# Row 51 (Stacks)
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.61"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  +9.07%  "

